$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.303.39'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').Value = '2.178.82'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '237.98'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.40'
$ws.Range('E7').Value = '  -3.14%  '
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.40'
$ws.Range('E10').Value = '  -5.19%  '
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').Value = '54.26'
$ws.Range('E12').Value = '  -5.65%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.101'
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '6.79'
$ws.Range('E14').Value = '  -3.54%  '
$ws.Range('D15').Value = '2.500.69'
$ws.Range('E15').Value = '  -1.35%  '
$ws.Range('D16').Value = '14.03'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '0.804'
$ws.Range('E17').Value = '  -3.78%  '
$ws.Range('D18').Value = '2.175.12'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = '41.091.83'
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('D20').Value = '0.0000102'
$ws.Range('E20').Value = '  -5.51%  '
$ws.Range('D21').Value = '70.55'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').Value = '5.96'
$ws.Range('E22').Value = '  -2.12%  '
$ws.Range('D23').Value = '9.85'
$ws.Range('E23').Value = '  -3.87%  '
$ws.Range('D24').Value = '226.42'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('E25').Value = '  -5.56%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '10.93'
$ws.Range('E27').Value = '  -4.45%  '
$ws.Range('D28').Value = '3.56'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('D31').Value = '167.86'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').Value = '19.99'
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('D33').Value = '31.46'
$ws.Range('E33').Value = '  +8.83%  '
$ws.Range('D34').Value = '0.0771'
$ws.Range('E34').Value = '  -2.05%  '
$ws.Range('D35').Value = '5.17'
$ws.Range('E35').Value = '  -6.65%  '
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('D37').Value = '0.104'
$ws.Range('E37').Value = '  -6.32%  '
$ws.Range('D38').Value = '4.13'
$ws.Range('E38').Value = '  -2.60%  '
$ws.Range('D39').Value = '0.0287'
$ws.Range('E39').Value = '  -3.97%  '
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('D41').Value = '11.88'
$ws.Range('E41').Value = '  -7.80%  '
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('D43').Value = '60.34'
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('D45').Value = '0.0979'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').Value = '8.31'
$ws.Range('E46').Value = '  -4.21%  '
$ws.Range('D47').Value = '98.53'
$ws.Range('E47').Value = '  -4.73%  '
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('E50').Value = '  -6.46%  '
$ws.Range('E51').Value = '  -2.66%  '
